# Natmi following Dr Hou advice
# Update the LR-pair results table (Gdf7 -> Bmpr2) on Sheet1:
#  - Rows 2-4 ("Sending cluster" = sCs) are recomputed with new stats and
#    the sending cluster relabelled to "FAPs".
#  - Three additional rows (5-7) are appended for the "sCs" sending
#    cluster against the same three target clusters (ECs, FAPs, sCs).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues {
    param($r, $vals)
    $c = 1
    foreach ($v in $vals) {
        $ws.Cells.Item($r, $c).Value = $v
        $c = $c + 1
    }
}

# Columns: A Sending cluster | B Ligand symbol | C Receptor symbol | D Target cluster
# E Ligand-expressing cells | F Ligand detection rate | G Ligand average expression value
# H Ligand total expression value | I Ligand derived specificity (avg) | J Ligand derived specificity (total)
# K Receptor-expressing cells | L Receptor detection rate | M Receptor average expression value
# N Receptor total expression value | O Receptor derived specificity (avg) | P Receptor derived specificity (total)
# Q Edge average expression weight | R Edge total expression weight
# S Edge average expression derived specificity | T Edge total expression derived specificity

Set-RowValues 2 @("FAPs","Gdf7","Bmpr2","ECs",2,0.6666666666666666,0.210033,0.630099,0.311901169594117,0.311901169594117,3,1,28.86405866666667,86.59217600000001,0.3025739760541936,0.3025739760541936,6.062404833936,54.561643505424,0.09437317702004532,0.09437317702004532)
Set-RowValues 3 @("FAPs","Gdf7","Bmpr2","FAPs",2,0.6666666666666666,0.210033,0.630099,0.311901169594117,0.311901169594117,3,1,42.68037399999999,128.041122,0.4474066038250392,0.4474066038250392,8.964286992341998,80.67858293107798,0.1395466430171614,0.1395466430171614)
Set-RowValues 4 @("FAPs","Gdf7","Bmpr2","sCs",2,0.6666666666666666,0.210033,0.630099,0.311901169594117,0.311901169594117,3,1,23.85061433333334,71.551843,0.2500194201207672,0.2500194201207672,5.009416080273001,45.084744722457,0.07798134955691019,0.07798134955691019)
Set-RowValues 5 @("sCs","Gdf7","Bmpr2","ECs",3,1,0.463363,1.390089,0.6880988304058829,0.6880988304058829,3,1,28.86405866666667,86.59217600000001,0.3025739760541936,0.3025739760541936,13.37453681596267,120.370831343664,0.2082007990341482,0.2082007990341482)
Set-RowValues 6 @("sCs","Gdf7","Bmpr2","FAPs",3,1,0.463363,1.390089,0.6880988304058829,0.6880988304058829,3,1,42.68037399999999,128.041122,0.4474066038250392,0.4474066038250392,19.776506137762,177.988555239858,0.3078599608078777,0.3078599608078777)
Set-RowValues 7 @("sCs","Gdf7","Bmpr2","sCs",3,1,0.463363,1.390089,0.6880988304058829,0.6880988304058829,3,1,23.85061433333334,71.551843,0.2500194201207672,0.2500194201207672,11.05149220933633,99.46342988402701,0.172038070563857,0.172038070563857)
